$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill the season record (Wins=51, Losses=64, Ties=0) for every data row (2..44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 29).Value = 51
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 0
}
